# Update sheet title and data to reflect new "through" date (2022-04-25)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the worksheet/tab to reflect the new date
$ws.Name = "Through 2022-04-25"

# Update the label in A5 ("April (through 04-23)" -> "April (through 04-25)")
$ws.Range("A5").Value = "April (through 04-25)"

# Update April row (row 5) values for years 2015-2022 (columns B-I)
$ws.Range("B5").Value = 16
$ws.Range("C5").Value = 29
$ws.Range("D5").Value = 50
$ws.Range("E5").Value = 42
$ws.Range("F5").Value = 37
$ws.Range("G5").Value = 52
$ws.Range("H5").Value = 88
$ws.Range("I5").Value = 106

# Update Total row (row 6) values for years 2015-2022 (columns B-I)
$ws.Range("B6").Value = 82
$ws.Range("C6").Value = 157
$ws.Range("D6").Value = 239
$ws.Range("E6").Value = 239
$ws.Range("F6").Value = 147
$ws.Range("G6").Value = 250
$ws.Range("H6").Value = 511
$ws.Range("I6").Value = 541
